$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 3's data, with new B value)
$ws.Range("A2").Value = 112182534
$ws.Range("B2").Value = 77402
$ws.Range("E2").Value = 6446
$ws.Range("F2").Value = "Kolflarnlav"
$ws.Range("G2").Value = "Carbonicola anthracophila"
$ws.Range("H2").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q2").Value = 431104
$ws.Range("R2").Value = 6811805

# Row 3 (was row 4's data, with new B value)
$ws.Range("A3").Value = 112181898
$ws.Range("B3").Value = 78216
$ws.Range("E3").Value = 229821
$ws.Range("F3").Value = "Vedflamlav"
$ws.Range("G3").Value = "Ramboldia elabens"
$ws.Range("H3").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q3").Value = 431104
$ws.Range("R3").Value = 6811804

# Row 4 (was row 2's data, with new B value)
$ws.Range("A4").Value = 112181853
$ws.Range("B4").Value = 78242
$ws.Range("E4").Value = 6453
$ws.Range("F4").Value = "Vedskivlav"
$ws.Range("G4").Value = "Hertelidea botryosa"
$ws.Range("H4").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q4").Value = 431106
$ws.Range("R4").Value = 6811802
